$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.374.73'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.21%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.881.40'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7129'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.08%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.46'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.08034'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +3.55%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.76%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08343'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -2.18%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.898.32'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.92%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.252'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.63%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.7190'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.17%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '94.13'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.89%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.339'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +5.50%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008563'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +4.35%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '29.390.64'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.28%  '
$ws.Range('B19').NumberFormat = '@'
$ws.Range('B19').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C19').NumberFormat = '@'
$ws.Range('C19').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.160.21'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.14%  '
$ws.Range('B20').NumberFormat = '@'
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').NumberFormat = '@'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '242.10'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.07%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.25'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.871'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.85%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1590'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -1.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '163.61'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.54%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.085'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.34%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.70%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.513'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.06%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.47%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.326'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.19%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -6.43%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +2.40%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.950'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.83%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.69%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7509'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.89%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.40%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01891'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.30%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.288.95'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +8.94%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.746'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.21%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +3.31%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9169'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +3.23%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '74.98'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +2.76%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '111.86'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +5.22%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.05%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000129'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +5.69%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.034.46'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.19%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.812'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.02%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.5222'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.29%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.528'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.50%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.83%  '
